{"js": "// Apply the edit described by the diff:\n//  1) \"...lots of different game; people \" -> \"...lots of different games; people \"\n//  2) \"...people who not like...\"          -> \"...people who do not like...\"\n//  3) \"...fewer game not having them...\"   -> \"...fewer games not having them...\"\n//\n// These are the only changes with a visible text impact; the remaining\n// hunks in the source diff only re-split/re-merge runs that already carry\n// identical formatting (no visible text difference between the before and\n// after runs), so there is nothing further to reproduce for those hunks.\n\nconst body = context.document.body;\n\n// --- Edit 1: \"different game;\" -> \"different games;\" -------------------\nconst gameResults = body.search(\"different game\", { matchCase: false, matchWholeWord: false });\ngameResults.load(\"text\");\nawait context.sync();\nif (gameResults.items.length !== 1) {\n  throw new Error('Expected exactly one match for \"different game\", found ' + gameResults.items.length);\n}\ngameResults.items[0].insertText(\"s\", Word.InsertLocation.end);\nawait context.sync();\n\n// --- Edit 2: \"people who not like\" -> \"people who do not like\" ---------\n// \"who not like\" is unique in the document; locate it, then narrow down\n// to the \"who\" token inside that match so \" do\" is inserted right after it.\nconst whoContext = body.search(\"who not like\", { matchCase: false, matchWholeWord: false });\nwhoContext.load(\"text\");\nawait context.sync();\nif (whoContext.items.length !== 1) {\n  throw new Error('Expected exactly one match for \"who not like\", found ' + whoContext.items.length);\n}\nconst whoToken = whoContext.items[0].search(\"who\", { matchCase: false, matchWholeWord: false });\nwhoToken.load(\"text\");\nawait context.sync();\nif (whoToken.items.length !== 1) {\n  throw new Error('Expected exactly one \"who\" inside \"who not like\" match.');\n}\nwhoToken.items[0].insertText(\" do\", Word.InsertLocation.end);\nawait context.sync();\n\n// --- Edit 3: \"fewer game not having them\" -> \"fewer games not having them\"\nconst gameContext2 = body.search(\"fewer game not having them\", { matchCase: false, matchWholeWord: false });\ngameContext2.load(\"text\");\nawait context.sync();\nif (gameContext2.items.length !== 1) {\n  throw new Error('Expected exactly one match for \"fewer game not having them\", found ' + gameContext2.items.length);\n}\nconst gameToken2 = gameContext2.items[0].search(\"game\", { matchCase: false, matchWholeWord: false });\ngameToken2.load(\"text\");\nawait context.sync();\nif (gameToken2.items.length !== 1) {\n  throw new Error('Expected exactly one \"game\" inside \"fewer game not having them\" match.');\n}\ngameToken2.items[0].insertText(\"s\", Word.InsertLocation.end);\nawait context.sync();\n", "ps1": "# Apply the edit described by the diff:\n#  1) \"...lots of different game; people \" -> \"...lots of different games; people \"\n#  2) \"...people who not like...\"          -> \"...people who do not like...\"\n#  3) \"...fewer game not having them...\"   -> \"...fewer games not having them...\"\n#\n# These are the only changes with a visible text impact; the remaining\n# hunks in the source diff only re-split/re-merge runs that already carry\n# identical formatting (no visible text difference between the before and\n# after runs), so there is nothing further to reproduce for those hunks.\n\n$d = $word.ActiveDocument\n\n# --- Edit 1: \"different game;\" -> \"different games;\" --------------------\n$r1 = $d.Content\n$found1 = $r1.Find.Execute(\"different game\")\nif (-not $found1) {\n    throw 'Could not find \"different game\" to pluralize.'\n}\n$r1.InsertAfter(\"s\")\n\n# --- Edit 2: \"people who not like\" -> \"people who do not like\" ----------\n# \"who not like\" is unique in the document; locate it, then narrow the\n# range down to just the \"who\" token so \" do\" lands right after it.\n$r2 = $d.Content\n$found2 = $r2.Find.Execute(\"who not like\")\nif (-not $found2) {\n    throw 'Could not find \"who not like\" to insert \"do\".'\n}\n$whoRange = $d.Range($r2.Start, $r2.Start + 3)\nif ($whoRange.Text -ne \"who\") {\n    throw 'Unexpected text when narrowing to \"who\" token: ' + $whoRange.Text\n}\n$whoRange.InsertAfter(\" do\")\n\n# --- Edit 3: \"fewer game not having them\" -> \"fewer games not having them\"\n$r3 = $d.Content\n$found3 = $r3.Find.Execute(\"fewer game not having them\")\nif (-not $found3) {\n    throw 'Could not find \"fewer game not having them\" to pluralize.'\n}\n$gameRange = $d.Range($r3.Start + 6, $r3.Start + 10)\nif ($gameRange.Text -ne \"game\") {\n    throw 'Unexpected text when narrowing to \"game\" token: ' + $gameRange.Text\n}\n$gameRange.InsertAfter(\"s\")\n"}
